$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The underlying dataset changed (extra data source rows read in), so the
# REF/NEW unique-id counts and all the dependent KPI counts/percentages in
# the summary sheet shift. All values in this sheet are stored as text,
# so force a text number format before writing each value to avoid Excel
# auto-converting numeric-looking strings into numbers/percentages.

$changes = @(
    @("D2",  "10"),
    @("D3",  "8"),
    @("D5",  "16"),
    @("D6",  "3"),
    @("F6",  "23.1%"),
    @("D7",  "1"),
    @("E7",  "10"),
    @("F7",  "10.0%"),
    @("E8",  "7"),
    @("F8",  "28.6%"),
    @("D9",  "5"),
    @("E9",  "7"),
    @("F9",  "71.4%"),
    @("G9",  "2,4,7,8,9"),
    @("E10", "10"),
    @("F10", "30.0%"),
    @("E11", "8"),
    @("F11", "12.5%"),
    @("E12", "7"),
    @("F12", "28.6%"),
    @("D13", "1"),
    @("E13", "10"),
    @("F13", "10.0%"),
    @("G13", "5"),
    @("E14", "8"),
    @("F14", "12.5%")
)

foreach ($change in $changes) {
    $cellRef = $change[0]
    $newValue = $change[1]
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $newValue
}
